$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "K" column (G) values with newly computed results
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 1
